# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.570.15"
$ws.Range("E2").Value = "  -13.21%  "
$ws.Range("D3").Value = "2.325.33"
$ws.Range("E3").Value = "  -19.88%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'447.35"
$ws.Range("E5").Value = "  -15.11%  "
$ws.Range("D6").Value = "'121.40"
$ws.Range("E6").Value = "  -15.05%  "
$ws.Range("E8").Value = "  -15.07%  "
$ws.Range("D9").Value = "2.315.03"
$ws.Range("E9").Value = "  -20.41%  "
$ws.Range("D10").Value = "'5.31"
$ws.Range("E10").Value = "  -11.00%  "
$ws.Range("D11").Value = "'0.0865"
$ws.Range("E11").Value = "  -19.17%  "
$ws.Range("D12").Value = "'0.301"
$ws.Range("E12").Value = "  -16.00%  "
$ws.Range("D13").Value = "'0.121"
$ws.Range("D14").Value = "52.606.53"
$ws.Range("E14").Value = "  -13.11%  "
$ws.Range("D15").Value = "'18.74"
$ws.Range("E15").Value = "  -17.07%  "
$ws.Range("D17").Value = "2.326.29"
$ws.Range("E17").Value = "  -20.04%  "
$ws.Range("D18").Value = "'3.94"
$ws.Range("E18").Value = "  -20.71%  "
$ws.Range("D19").Value = "'297.43"
$ws.Range("E19").Value = "  -15.55%  "
$ws.Range("D20").Value = "'8.88"
$ws.Range("E20").Value = "  -23.59%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "'5.63"
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("D23").Value = "'5.13"
$ws.Range("E23").Value = "  -21.73%  "
$ws.Range("D24").Value = "'53.51"
$ws.Range("E24").Value = "  -17.20%  "
$ws.Range("E25").Value = "  -19.48%  "
$ws.Range("D26").Value = "'0.145"
$ws.Range("E26").Value = "  -18.55%  "
$ws.Range("E27").Value = "  -11.60%  "
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("D29").Value = "0.0₃0659"
$ws.Range("E29").Value = "  -22.25%  "
$ws.Range("D30").Value = "'140.17"
$ws.Range("E30").Value = "  -7.24%  "
$ws.Range("D31").Value = "'16.80"
$ws.Range("E31").Value = "  -14.25%  "
$ws.Range("E32").Value = "  -19.72%  "
$ws.Range("E33").Value = "  -15.37%  "
$ws.Range("D34").Value = "'0.821"
$ws.Range("E34").Value = "  -17.73%  "
$ws.Range("E35").Value = "  -21.52%  "
$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "'0.992"
$ws.Range("E37").Value = "  -17.04%  "
$ws.Range("D38").Value = "'31.79"
$ws.Range("E38").Value = "  -15.66%  "
$ws.Range("D39").Value = "'10.14"
$ws.Range("E39").Value = "  -1.77%  "
$ws.Range("D40").Value = "'0.560"
$ws.Range("E40").Value = "  -13.63%  "
$ws.Range("D41").Value = "'0.0505"
$ws.Range("E41").Value = "  -13.14%  "
$ws.Range("D42").Value = "'3.10"
$ws.Range("E42").Value = "  -16.53%  "
$ws.Range("D43").Value = "1.903.52"
$ws.Range("E43").Value = "  -16.92%  "
$ws.Range("D44").Value = "'1.17"
$ws.Range("E44").Value = "  -20.47%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0205"
$ws.Range("E45").Value = "  -13.54%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.0820"
$ws.Range("E46").Value = "  -10.86%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'4.18"
$ws.Range("E47").Value = "  -15.39%  "
$ws.Range("D48").Value = "'15.59"
$ws.Range("E48").Value = "  -23.59%  "
$ws.Range("D49").Value = "'4.59"
$ws.Range("E49").Value = "  -5.18%  "
$ws.Range("D50").Value = "'4.44"
$ws.Range("E50").Value = "  -12.98%  "
$ws.Range("D51").Value = "'15.01"
$ws.Range("E51").Value = "  -17.77%  "
